# This script adds 20 new ranked rows (rows 52-71) to each of the 5
# worksheets in the workbook, extending the word-count ranking lists
# that previously stopped at row 51 (dimension A1:C51 -> A1:C71).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "전체" ---
$ws = $wb.Worksheets.Item("전체")
$ws.Cells.Item(52, 1).Value = 484
$ws.Cells.Item(52, 2).Value = "배당 "
$ws.Cells.Item(52, 3).Value = 31399
$ws.Cells.Item(53, 1).Value = 1509
$ws.Cells.Item(53, 2).Value = "행사 가격"
$ws.Cells.Item(53, 3).Value = 30943
$ws.Cells.Item(54, 1).Value = 327
$ws.Cells.Item(54, 2).Value = "환율"
$ws.Cells.Item(54, 3).Value = 30930
$ws.Cells.Item(55, 1).Value = 1051
$ws.Cells.Item(55, 2).Value = "그린채권"
$ws.Cells.Item(55, 3).Value = 30483
$ws.Cells.Item(56, 1).Value = 1744
$ws.Cells.Item(56, 2).Value = "중화"
$ws.Cells.Item(56, 3).Value = 29320
$ws.Cells.Item(57, 1).Value = 586
$ws.Cells.Item(57, 2).Value = "유럽연합"
$ws.Cells.Item(57, 3).Value = 29071
$ws.Cells.Item(58, 1).Value = 1013
$ws.Cells.Item(58, 2).Value = "결산"
$ws.Cells.Item(58, 3).Value = 28483
$ws.Cells.Item(59, 1).Value = 1240
$ws.Cells.Item(59, 2).Value = "시가총액"
$ws.Cells.Item(59, 3).Value = 28475
$ws.Cells.Item(60, 1).Value = 1104
$ws.Cells.Item(60, 2).Value = "디지털 금융"
$ws.Cells.Item(60, 3).Value = 28419
$ws.Cells.Item(61, 1).Value = 1698
$ws.Cells.Item(61, 2).Value = "입찰"
$ws.Cells.Item(61, 3).Value = 27618
$ws.Cells.Item(62, 1).Value = 1450
$ws.Cells.Item(62, 2).Value = "부가가치"
$ws.Cells.Item(62, 3).Value = 27587
$ws.Cells.Item(63, 1).Value = 541
$ws.Cells.Item(63, 2).Value = "시공사"
$ws.Cells.Item(63, 3).Value = 26038
$ws.Cells.Item(64, 1).Value = 686
$ws.Cells.Item(64, 2).Value = "추가경정예산"
$ws.Cells.Item(64, 3).Value = 25981
$ws.Cells.Item(65, 1).Value = 423
$ws.Cells.Item(65, 2).Value = "국세"
$ws.Cells.Item(65, 3).Value = 25943
$ws.Cells.Item(66, 1).Value = 2025
$ws.Cells.Item(66, 2).Value = "독점"
$ws.Cells.Item(66, 3).Value = 25828
$ws.Cells.Item(67, 1).Value = 613
$ws.Cells.Item(67, 2).Value = "재건축"
$ws.Cells.Item(67, 3).Value = 25610
$ws.Cells.Item(68, 1).Value = 833
$ws.Cells.Item(68, 2).Value = "소셜네트워크서비스"
$ws.Cells.Item(68, 3).Value = 25303
$ws.Cells.Item(69, 1).Value = 187
$ws.Cells.Item(69, 2).Value = "작은 정부"
$ws.Cells.Item(69, 3).Value = 24791
$ws.Cells.Item(70, 1).Value = 944
$ws.Cells.Item(70, 2).Value = "당기순이익"
$ws.Cells.Item(70, 3).Value = 24115
$ws.Cells.Item(71, 1).Value = 1704
$ws.Cells.Item(71, 2).Value = "자기자본"
$ws.Cells.Item(71, 3).Value = 23921
# Copy the bold/bordered style used by column A (e.g. A51) down
# through the newly added rows so formatting matches the rest of
# the ranking table.
$ws.Range("A51").Copy()
$ws.Range("A52:A71").PasteSpecial(-4122)

# --- Sheet 2: "사회" ---
$ws = $wb.Worksheets.Item("사회")
$ws.Cells.Item(52, 1).Value = 398
$ws.Cells.Item(52, 2).Value = "재개발"
$ws.Cells.Item(52, 3).Value = 4363
$ws.Cells.Item(53, 1).Value = 320
$ws.Cells.Item(53, 2).Value = "부담금"
$ws.Cells.Item(53, 3).Value = 4358
$ws.Cells.Item(54, 1).Value = 6
$ws.Cells.Item(54, 2).Value = "범위의 경제"
$ws.Cells.Item(54, 3).Value = 4358
$ws.Cells.Item(55, 1).Value = 769
$ws.Cells.Item(55, 2).Value = "빅뱅"
$ws.Cells.Item(55, 3).Value = 4346
$ws.Cells.Item(56, 1).Value = 1030
$ws.Cells.Item(56, 2).Value = "수재"
$ws.Cells.Item(56, 3).Value = 4308
$ws.Cells.Item(57, 1).Value = 518
$ws.Cells.Item(57, 2).Value = "디지털 경제"
$ws.Cells.Item(57, 3).Value = 4304
$ws.Cells.Item(58, 1).Value = 1038
$ws.Cells.Item(58, 2).Value = "CBS"
$ws.Cells.Item(58, 3).Value = 4285
$ws.Cells.Item(59, 1).Value = 105
$ws.Cells.Item(59, 2).Value = "인구 고령화"
$ws.Cells.Item(59, 3).Value = 4258
$ws.Cells.Item(60, 1).Value = 695
$ws.Cells.Item(60, 2).Value = "금융감독원"
$ws.Cells.Item(60, 3).Value = 4164
$ws.Cells.Item(61, 1).Value = 1274
$ws.Cells.Item(61, 2).Value = "기본 소득"
$ws.Cells.Item(61, 3).Value = 4152
$ws.Cells.Item(62, 1).Value = 1306
$ws.Cells.Item(62, 2).Value = "독점"
$ws.Cells.Item(62, 3).Value = 4070
$ws.Cells.Item(63, 1).Value = 404
$ws.Cells.Item(63, 2).Value = "저작권"
$ws.Cells.Item(63, 3).Value = 4068
$ws.Cells.Item(64, 1).Value = 593
$ws.Cells.Item(64, 2).Value = "트위터"
$ws.Cells.Item(64, 3).Value = 3974
$ws.Cells.Item(65, 1).Value = 965
$ws.Cells.Item(65, 2).Value = "홀딩스"
$ws.Cells.Item(65, 3).Value = 3947
$ws.Cells.Item(66, 1).Value = 1245
$ws.Cells.Item(66, 2).Value = "공정거래위원회"
$ws.Cells.Item(66, 3).Value = 3916
$ws.Cells.Item(67, 1).Value = 781
$ws.Cells.Item(67, 2).Value = "상환"
$ws.Cells.Item(67, 3).Value = 3882
$ws.Cells.Item(68, 1).Value = 1207
$ws.Cells.Item(68, 2).Value = "PER"
$ws.Cells.Item(68, 3).Value = 3851
$ws.Cells.Item(69, 1).Value = 1040
$ws.Cells.Item(69, 2).Value = "EaR"
$ws.Cells.Item(69, 3).Value = 3742
$ws.Cells.Item(70, 1).Value = 146
$ws.Cells.Item(70, 2).Value = "직접 금융 "
$ws.Cells.Item(70, 3).Value = 3724
$ws.Cells.Item(71, 1).Value = 535
$ws.Cells.Item(71, 2).Value = "CEO"
$ws.Cells.Item(71, 3).Value = 3718
# Copy the bold/bordered style used by column A (e.g. A51) down
# through the newly added rows so formatting matches the rest of
# the ranking table.
$ws.Range("A51").Copy()
$ws.Range("A52:A71").PasteSpecial(-4122)

# --- Sheet 3: "경제" ---
$ws = $wb.Worksheets.Item("경제")
$ws.Cells.Item(52, 1).Value = 1383
$ws.Cells.Item(52, 2).Value = "홀딩스"
$ws.Cells.Item(52, 3).Value = 16305
$ws.Cells.Item(53, 1).Value = 1796
$ws.Cells.Item(53, 2).Value = "기준금리"
$ws.Cells.Item(53, 3).Value = 15725
$ws.Cells.Item(54, 1).Value = 1320
$ws.Cells.Item(54, 2).Value = "부가가치"
$ws.Cells.Item(54, 3).Value = 15518
$ws.Cells.Item(55, 1).Value = 1038
$ws.Cells.Item(55, 2).Value = "보합세"
$ws.Cells.Item(55, 3).Value = 14826
$ws.Cells.Item(56, 1).Value = 1193
$ws.Cells.Item(56, 2).Value = "유가증권"
$ws.Cells.Item(56, 3).Value = 14793
$ws.Cells.Item(57, 1).Value = 750
$ws.Cells.Item(57, 2).Value = "스마트 자동차"
$ws.Cells.Item(57, 3).Value = 14761
$ws.Cells.Item(58, 1).Value = 613
$ws.Cells.Item(58, 2).Value = "출연"
$ws.Cells.Item(58, 3).Value = 14596
$ws.Cells.Item(59, 1).Value = 1072
$ws.Cells.Item(59, 2).Value = "상품권"
$ws.Cells.Item(59, 3).Value = 14536
$ws.Cells.Item(60, 1).Value = 1743
$ws.Cells.Item(60, 2).Value = "고용없는 성장"
$ws.Cells.Item(60, 3).Value = 14502
$ws.Cells.Item(61, 1).Value = 1459
$ws.Cells.Item(61, 2).Value = "CTO"
$ws.Cells.Item(61, 3).Value = 14276
$ws.Cells.Item(62, 1).Value = 302
$ws.Cells.Item(62, 2).Value = "구조조정"
$ws.Cells.Item(62, 3).Value = 14240
$ws.Cells.Item(63, 1).Value = 1546
$ws.Cells.Item(63, 2).Value = "자본 경영"
$ws.Cells.Item(63, 3).Value = 13955
$ws.Cells.Item(64, 1).Value = 1772
$ws.Cells.Item(64, 2).Value = "국내총생산"
$ws.Cells.Item(64, 3).Value = 13954
$ws.Cells.Item(65, 1).Value = 822
$ws.Cells.Item(65, 2).Value = "CB"
$ws.Cells.Item(65, 3).Value = 13853
$ws.Cells.Item(66, 1).Value = 1730
$ws.Cells.Item(66, 2).Value = "경제성장률"
$ws.Cells.Item(66, 3).Value = 13764
$ws.Cells.Item(67, 1).Value = 1220
$ws.Cells.Item(67, 2).Value = "임대료"
$ws.Cells.Item(67, 3).Value = 13669
$ws.Cells.Item(68, 1).Value = 184
$ws.Cells.Item(68, 2).Value = "제품 차별화"
$ws.Cells.Item(68, 3).Value = 13200
$ws.Cells.Item(69, 1).Value = 1631
$ws.Cells.Item(69, 2).Value = "판매 채널"
$ws.Cells.Item(69, 3).Value = 13167
$ws.Cells.Item(70, 1).Value = 526
$ws.Cells.Item(70, 2).Value = "유럽연합"
$ws.Cells.Item(70, 3).Value = 13106
$ws.Cells.Item(71, 1).Value = 1352
$ws.Cells.Item(71, 2).Value = "핀테크"
$ws.Cells.Item(71, 3).Value = 12881
# Copy the bold/bordered style used by column A (e.g. A51) down
# through the newly added rows so formatting matches the rest of
# the ranking table.
$ws.Range("A51").Copy()
$ws.Range("A52:A71").PasteSpecial(-4122)

# --- Sheet 4: "문화" ---
$ws = $wb.Worksheets.Item("문화")
$ws.Cells.Item(52, 1).Value = 344
$ws.Cells.Item(52, 2).Value = "재개발"
$ws.Cells.Item(52, 3).Value = 2311
$ws.Cells.Item(53, 1).Value = 10
$ws.Cells.Item(53, 2).Value = "부의 효과"
$ws.Cells.Item(53, 3).Value = 2257
$ws.Cells.Item(54, 1).Value = 648
$ws.Cells.Item(54, 2).Value = "빅뱅"
$ws.Cells.Item(54, 3).Value = 2252
$ws.Cells.Item(55, 1).Value = 535
$ws.Cells.Item(55, 2).Value = "LaR"
$ws.Cells.Item(55, 3).Value = 2131
$ws.Cells.Item(56, 1).Value = 459
$ws.Cells.Item(56, 2).Value = "스마트 교육"
$ws.Cells.Item(56, 3).Value = 2014
$ws.Cells.Item(57, 1).Value = 487
$ws.Cells.Item(57, 2).Value = "증강현실"
$ws.Cells.Item(57, 3).Value = 1939
$ws.Cells.Item(58, 1).Value = 486
$ws.Cells.Item(58, 2).Value = "점포자동화"
$ws.Cells.Item(58, 3).Value = 1891
$ws.Cells.Item(59, 1).Value = 663
$ws.Cells.Item(59, 2).Value = "소득 대비 대출 비율"
$ws.Cells.Item(59, 3).Value = 1884
$ws.Cells.Item(60, 1).Value = 605
$ws.Cells.Item(60, 2).Value = "넉아웃"
$ws.Cells.Item(60, 3).Value = 1883
$ws.Cells.Item(61, 1).Value = 447
$ws.Cells.Item(61, 2).Value = "빅데이터"
$ws.Cells.Item(61, 3).Value = 1868
$ws.Cells.Item(62, 1).Value = 1100
$ws.Cells.Item(62, 2).Value = "그린슈트"
$ws.Cells.Item(62, 3).Value = 1863
$ws.Cells.Item(63, 1).Value = 274
$ws.Cells.Item(63, 2).Value = "배당 "
$ws.Cells.Item(63, 3).Value = 1803
$ws.Cells.Item(64, 1).Value = 1098
$ws.Cells.Item(64, 2).Value = "규모의 경제"
$ws.Cells.Item(64, 3).Value = 1764
$ws.Cells.Item(65, 1).Value = 899
$ws.Cells.Item(65, 2).Value = "엑시트"
$ws.Cells.Item(65, 3).Value = 1727
$ws.Cells.Item(66, 1).Value = 580
$ws.Cells.Item(66, 2).Value = "COO"
$ws.Cells.Item(66, 3).Value = 1705
$ws.Cells.Item(67, 1).Value = 842
$ws.Cells.Item(67, 2).Value = "밀레니얼 세대"
$ws.Cells.Item(67, 3).Value = 1657
$ws.Cells.Item(68, 1).Value = 1003
$ws.Cells.Item(68, 2).Value = "ABC"
$ws.Cells.Item(68, 3).Value = 1569
$ws.Cells.Item(69, 1).Value = 779
$ws.Cells.Item(69, 2).Value = "콜"
$ws.Cells.Item(69, 3).Value = 1551
$ws.Cells.Item(70, 1).Value = 700
$ws.Cells.Item(70, 2).Value = "연금"
$ws.Cells.Item(70, 3).Value = 1523
$ws.Cells.Item(71, 1).Value = 185
$ws.Cells.Item(71, 2).Value = "환경 경영"
$ws.Cells.Item(71, 3).Value = 1520
# Copy the bold/bordered style used by column A (e.g. A51) down
# through the newly added rows so formatting matches the rest of
# the ranking table.
$ws.Range("A51").Copy()
$ws.Range("A52:A71").PasteSpecial(-4122)

# --- Sheet 5: "IT" ---
$ws = $wb.Worksheets.Item("IT")
$ws.Cells.Item(52, 1).Value = 1311
$ws.Cells.Item(52, 2).Value = "규모의 경제"
$ws.Cells.Item(52, 3).Value = 5109
$ws.Cells.Item(53, 1).Value = 700
$ws.Cells.Item(53, 2).Value = "금융감독원"
$ws.Cells.Item(53, 3).Value = 5053
$ws.Cells.Item(54, 1).Value = 631
$ws.Cells.Item(54, 2).Value = "GP"
$ws.Cells.Item(54, 3).Value = 5035
$ws.Cells.Item(55, 1).Value = 598
$ws.Cells.Item(55, 2).Value = "플레이 스토어"
$ws.Cells.Item(55, 3).Value = 4961
$ws.Cells.Item(56, 1).Value = 31
$ws.Cells.Item(56, 2).Value = "선별"
$ws.Cells.Item(56, 3).Value = 4725
$ws.Cells.Item(57, 1).Value = 544
$ws.Cells.Item(57, 2).Value = "스마트 팩토리"
$ws.Cells.Item(57, 3).Value = 4628
$ws.Cells.Item(58, 1).Value = 995
$ws.Cells.Item(58, 2).Value = "행사 가격"
$ws.Cells.Item(58, 3).Value = 4552
$ws.Cells.Item(59, 1).Value = 596
$ws.Cells.Item(59, 2).Value = "트위터"
$ws.Cells.Item(59, 3).Value = 4515
$ws.Cells.Item(60, 1).Value = 958
$ws.Cells.Item(60, 2).Value = "부가가치"
$ws.Cells.Item(60, 3).Value = 4444
$ws.Cells.Item(61, 1).Value = 638
$ws.Cells.Item(61, 2).Value = "MOR"
$ws.Cells.Item(61, 3).Value = 4305
$ws.Cells.Item(62, 1).Value = 454
$ws.Cells.Item(62, 2).Value = "4G"
$ws.Cells.Item(62, 3).Value = 4229
$ws.Cells.Item(63, 1).Value = 469
$ws.Cells.Item(63, 2).Value = "OLED"
$ws.Cells.Item(63, 3).Value = 4062
$ws.Cells.Item(64, 1).Value = 12
$ws.Cells.Item(64, 2).Value = "부의 효과"
$ws.Cells.Item(64, 3).Value = 4043
$ws.Cells.Item(65, 1).Value = 615
$ws.Cells.Item(65, 2).Value = "CB"
$ws.Cells.Item(65, 3).Value = 4012
$ws.Cells.Item(66, 1).Value = 795
$ws.Cells.Item(66, 2).Value = "상품권"
$ws.Cells.Item(66, 3).Value = 3988
$ws.Cells.Item(67, 1).Value = 946
$ws.Cells.Item(67, 2).Value = "컨소시엄"
$ws.Cells.Item(67, 3).Value = 3984
$ws.Cells.Item(68, 1).Value = 553
$ws.Cells.Item(68, 2).Value = "와이파이"
$ws.Cells.Item(68, 3).Value = 3886
$ws.Cells.Item(69, 1).Value = 352
$ws.Cells.Item(69, 2).Value = "양해각서"
$ws.Cells.Item(69, 3).Value = 3804
$ws.Cells.Item(70, 1).Value = 277
$ws.Cells.Item(70, 2).Value = "기금"
$ws.Cells.Item(70, 3).Value = 3715
$ws.Cells.Item(71, 1).Value = 1294
$ws.Cells.Item(71, 2).Value = "공정거래위원회"
$ws.Cells.Item(71, 3).Value = 3578
# Copy the bold/bordered style used by column A (e.g. A51) down
# through the newly added rows so formatting matches the rest of
# the ranking table.
$ws.Range("A51").Copy()
$ws.Range("A52:A71").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Host "Done: extended all 5 sheets from 51 to 71 rows."
